$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.760.41'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '1.628.12'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.57'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = '1.853.80'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').Value = '1.630.15'
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '0.0₃0755'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '25.778.66'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.997'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.997'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.30%  '
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.22'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.902'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').Value = '1.141.37'
$ws.Range('E37').Value = '  +2.06%  '
$ws.Range('E38').Value = '  -0.68%  '
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.996'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.62'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.799'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.33%  '
$ws.Range('D45').Value = '1.765.24'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0512'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.17%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.52%  '
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0954'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.78%  '
